$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "вадик"
$ws.Range("B5").Value = "CAACAgIAAxkBAAICB2I5jl5j1LKESCtSIHLFgs05kXD2AAIoEgACsJtxSHvFqJ5Z81PdIwQ"
